# Scene.xlsx - add "CanClone" and "ActorID" columns to the NPC scene table,
# fill in their values for the 3 existing rows, widen the new "CanClone"
# column, flip the page to portrait / Letter-ish (paper size 9), and leave
# the selection on K9 (matches the author's edit captured in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Append two new columns to the XML-mapped table. ListColumns.Add() always
# appends at the next free position, so the first call lands in J, the
# second in K.
$colJ = $lo.ListColumns.Add()
$colK = $lo.ListColumns.Add()

# Name the headers. Do K first, then J: shared-string ids are handed out in
# write order, and the target workbook has "ActorID" interned before
# "CanClone" even though CanClone is the left-hand (J) column.
$ws.Range("K1").Value2 = "ActorID"
$ws.Range("J1").Value2 = "CanClone"

# Row data for the 3 existing records.
$ws.Range("J2").Value2 = 1
$ws.Range("K2").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 0

# New "CanClone" column gets an explicit width (14 chars stored).
$ws.Columns.Item(10).ColumnWidth = 13.29

# Page setup: portrait, paper size 9 (A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection left on K9 by the author.
[void]$ws.Range("K9").Select()
